# Weekly data refresh: a new "Terminal La Palmera de La Serena" Acelga
# price-report date pair is prepended to the historical table (rows 168-169),
# pushing every existing date pair down by one position. The oldest pair that
# falls off the bottom of its old slot re-appears duplicated in the two
# brand-new rows created at the end of the sheet (260-261).
#
# Implemented as: insert two blank rows at the top of the block (which Excel
# shifts the whole rest of the table down by two, automatically carrying the
# date-style formatting onto the new blank rows), seed those two rows with a
# copy of what is now immediately below them (the original first pair), then
# overwrite just the Fecha (date) and Volumen (J) values for the new pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole historical block (old rows 168:259) down by two rows.
$ws.Rows.Item(168).Insert()
$ws.Rows.Item(168).Insert()

# Seed the two new rows with the data now sitting right below them (this was
# the original first pair, now at 170:171), picking up every column
# (region/category/quality/price-range/unit/origin/etc.) and the date style.
$ws.Range("A170:R171").Copy($ws.Range("A168:R169"))

# Overwrite with the real new pair's date and volume.
$ws.Cells.Item(168, 4).Value = 44510
$ws.Cells.Item(168, 10).Value = 2760

$ws.Cells.Item(169, 4).Value = 44510
$ws.Cells.Item(169, 10).Value = 1400
